$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Value
    )
    $rng = $ws.Range($Address)
    $rng.NumberFormat = "@"
    $rng.Value = $Value
    $rng.Style = "Normal"
}

# Row 41 previously had an (empty) placeholder cell in column D.
# The updated log no longer has any value there, so remove it entirely.
$ws.Range("D41").ClearContents()

# New log rows appended on 25/03/2025

# Row 42
Set-TextValue "A42" "25/03/2025"
Set-TextValue "B42" "18:48:12"
Set-TextValue "C42" "invoice.pdf"
Set-TextValue "D42" "100"
Set-TextValue "E42" "COMPLETED"
Set-TextValue "F42" "OK"

# Row 43
Set-TextValue "A43" "25/03/2025"
Set-TextValue "B43" "18:48:13"
Set-TextValue "C43" "Invoice1.pdf"
Set-TextValue "E43" "Exception"
Set-TextValue "F43" "Invoice outside the specified standard: 'NoneType' object has no attribute 'group'"

# Row 44
Set-TextValue "A44" "25/03/2025"
Set-TextValue "B44" "18:49:16"
Set-TextValue "C44" "invoice.pdf"
Set-TextValue "D44" "100"
Set-TextValue "E44" "COMPLETED"
Set-TextValue "F44" "OK"

# Row 45
Set-TextValue "A45" "25/03/2025"
Set-TextValue "B45" "18:49:16"
Set-TextValue "C45" "Invoice1.pdf"
Set-TextValue "E45" "Exception"
Set-TextValue "F45" "Invoice outside the specified standard: 'NoneType' object has no attribute 'group'"

# Row 46
Set-TextValue "A46" "25/03/2025"
Set-TextValue "B46" "19:31:24"
Set-TextValue "C46" "invoice.pdf"
Set-TextValue "D46" "100"
Set-TextValue "E46" "COMPLETED"
Set-TextValue "F46" "OK"

# Row 47
Set-TextValue "A47" "25/03/2025"
Set-TextValue "B47" "19:31:25"
Set-TextValue "C47" "Invoice1.pdf"
# D47 stays blank (matches the blank placeholder seen on the previous
# exception row), so nothing is written there.
Set-TextValue "E47" "Exception"
Set-TextValue "F47" "Invoice outside the specified standard: 'NoneType' object has no attribute 'group'"
